$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F = 想去人数 (want-to-go count), column G = 最低票价 (lowest price)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 607
$ws1.Range("F5").Value = 1143
$ws1.Range("F6").Value = 14246
$ws1.Range("F7").Value = 16156
$ws1.Range("G7").Value = 65
$ws1.Range("F20").Value = 1238
$ws1.Range("F23").Value = 28
$ws1.Range("F24").Value = 6472
$ws1.Range("F26").Value = 8
$ws1.Range("F29").Value = 5674
$ws1.Range("F30").Value = 92
$ws1.Range("F33").Value = 4704
$ws1.Range("F34").Value = 13

# Sheet "全部类型" (sheet4) - same updates, rows shifted due to extra entries in this sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 607
$ws4.Range("F5").Value = 1143
$ws4.Range("F6").Value = 14246
$ws4.Range("F7").Value = 16156
$ws4.Range("G7").Value = 65
$ws4.Range("F20").Value = 1238
$ws4.Range("F24").Value = 28
$ws4.Range("F25").Value = 6472
$ws4.Range("F27").Value = 8
$ws4.Range("F31").Value = 5674
$ws4.Range("F32").Value = 92
$ws4.Range("F35").Value = 4704
$ws4.Range("F36").Value = 13

$wb.Save()
